$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "List1"

# Clear row 10 (B10 had a "Zarez"/comma style applied; remove it with the row)
$ws.Rows.Item(10).Delete()

# Remove the custom number-format style from column B (B1:B9) and set new values
$ws.Range("B1:B9").ClearFormats()

# The "Zarez" (comma) cell style is now unused; remove its definition too
$wb.Styles.Item("Zarez").Delete()

# New data values (no formulas) for A1:B9
$values = @(
    @(1, 3),
    @(3, 4),
    @(7, 2),
    @(13, 2),
    @(17, 2),
    @(22, 2),
    @(28, 2),
    @(32, 2),
    @(43, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Update selection to match target (active cell A2)
$ws.Range("A2").Select()
